$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value while forcing it to remain a text string.
# This matters because several "Price" column entries look like numbers
# (e.g. "1.00", "65.692.78", "0.130") but must stay as text, matching
# the original inlineStr cell type, preserving formatting such as
# trailing zeros and "."-as-thousands-separator groupings.
function Set-TextValue($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '65.692.78'
Set-TextValue "D3" '3.749.07'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue "D5" '423.71'
$ws.Range("E5").Value = '  +4.89%  '
Set-TextValue "D6" '136.48'
$ws.Range("E6").Value = '  +6.15%  '
Set-TextValue "D7" '0.615'
$ws.Range("E7").Value = '  +2.12%  '
Set-TextValue "D8" '1.00'
$ws.Range("E8").Value = '  -0.06%  '
Set-TextValue "D9" '0.722'
$ws.Range("E9").Value = '  +0.76%  '
Set-TextValue "D10" '0.147'
$ws.Range("E10").Value = '  -10.78%  '
Set-TextValue "D11" '0.0000296'
$ws.Range("E11").Value = '  -16.69%  '
Set-TextValue "D12" '41.77'
$ws.Range("E12").Value = '  +3.79%  '
Set-TextValue "D13" '10.32'
$ws.Range("E13").Value = '  +7.60%  '
Set-TextValue "D14" '4.357.40'
$ws.Range("E14").Value = '  +0.56%  '
Set-TextValue "D15" '14.91'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("E16").Value = '  +0.20%  '
Set-TextValue "D17" '3.757.24'
$ws.Range("E17").Value = '  +0.44%  '
Set-TextValue "D18" '19.64'
$ws.Range("E18").Value = '  +1.57%  '
Set-TextValue "D19" '1.10'
$ws.Range("E19").Value = '  +4.92%  '
Set-TextValue "D20" '65.909.85'
$ws.Range("E20").Value = '  -0.69%  '
Set-TextValue "D21" '398.75'
$ws.Range("E21").Value = '  -1.92%  '
Set-TextValue "D22" '14.76'
$ws.Range("E22").Value = '  +2.81%  '
Set-TextValue "D23" '3.18'
$ws.Range("E23").Value = '  +6.08%  '
Set-TextValue "D24" '83.65'
$ws.Range("E24").Value = '  -1.35%  '
Set-TextValue "D25" '36.16'
$ws.Range("E25").Value = '  +0.86%  '
Set-TextValue "D26" '9.77'
$ws.Range("E26").Value = '  +37.99%  '
Set-TextValue "D27" '3.22'
$ws.Range("E27").Value = '  +4.63%  '
Set-TextValue "D28" '9.70'
$ws.Range("E28").Value = '  +4.30%  '
$ws.Range("E29").Value = '  -4.05%  '
Set-TextValue "D30" '13.57'
$ws.Range("E30").Value = '  +10.31%  '
Set-TextValue "D31" '699.59'
$ws.Range("E31").Value = '  +1.85%  '
Set-TextValue "D32" '0.130'
$ws.Range("E32").Value = '  +12.47%  '
Set-TextValue "D33" '2.75'
$ws.Range("E33").Value = '  +1.84%  '
Set-TextValue "D34" '40.25'
$ws.Range("E34").Value = '  +4.09%  '
Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  -0.01%  '
Set-TextValue "D36" '5.56'
$ws.Range("E36").Value = '  +32.10%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D37" '0.147'
$ws.Range("E37").Value = '  -4.68%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D38" '56.26'
$ws.Range("E38").Value = '  +2.25%  '
Set-TextValue "D39" '0.0464'
$ws.Range("E39").Value = '  +2.52%  '
Set-TextValue "D40" '2.60'
$ws.Range("E40").Value = '  +39.20%  '
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("E42").Value = '  +0.21%  '
Set-TextValue "D43" '0.139'
$ws.Range("E43").Value = '  +3.69%  '
Set-TextValue "D44" '0.0₃0642'
$ws.Range("E44").Value = '  -10.71%  '
$ws.Range("E45").Value = '  +1.29%  '
Set-TextValue "D46" '3.31'
$ws.Range("E46").Value = '  +3.25%  '
Set-TextValue "D47" '0.315'
$ws.Range("E47").Value = '  +9.34%  '
$ws.Range("E48").Value = '  +4.55%  '
Set-TextValue "D49" '2.02'
$ws.Range("E49").Value = '  -0.74%  '
Set-TextValue "D50" '140.32'
$ws.Range("E50").Value = '  -3.18%  '
Set-TextValue "D51" '2.73'
$ws.Range("E51").Value = '  -2.63%  '
